$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the three new "search items" rows below the existing list.
$ws.Range("A8").Value = "graphic card"
$ws.Range("A9").Value = "bike chain cleaner"
$ws.Range("A10").Value = "bike chain lube"

# Column B was only ever an empty placeholder next to the items list;
# clear it out (not a structural column delete - the column width stays).
$ws.Range("B1:B5").ClearContents()

# Leave the selection on the last entry, matching the saved view state.
$ws.Range("A10").Select()
